# refactor: change to PSA naming
# Renames the 3-letter sector codes used both as column headers (row 1)
# and as row labels (column A) to the new PSA-style codes, and adds a
# bottom border under the EXT / TOTAL rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rename = @{
  "AGR" = "AFF"
  "MIN" = "MAQ"
  "MAN" = "MFG"
  "ESW" = "ESWW"
  "CON" = "CNS"
  "WRT" = "TRD"
  "TRS" = "TAS"
  "AFS" = "AFSA"
  "INF" = "IAC"
  "FIN" = "FIA"
  "REA" = "REOD"
  "PBS" = "PBS"
  "PAD" = "PAD"
  "EDU" = "EDUC"
  "HHS" = "HHSW"
  "OTH" = "OS"
  "CAP" = "CAP"
  "LAB" = "LAB"
  "IDT" = "IDT"
  "TRF" = "TRF"
  "HOH" = "HOH"
  "GOV" = "GOV"
  "INV" = "INV"
  "EXT" = "EXT"
  "TOTAL" = "TOTAL"
}

# Header row, B1:Z1 - column labels
for ($c = 2; $c -le 26; $c++) {
    $old = $ws.Cells.Item(1, $c).Text
    if ($rename.ContainsKey($old)) {
        $ws.Cells.Item(1, $c).Value = $rename[$old]
    }
}

# Column A, rows 2:26 - row labels
for ($r = 2; $r -le 26; $r++) {
    $old = $ws.Cells.Item($r, 1).Text
    if ($rename.ContainsKey($old)) {
        $ws.Cells.Item($r, 1).Value = $rename[$old]
    }
}

# The (renamed) EXT row (A25) and TOTAL row (A26) labels now get a
# left+right+bottom thin border (previously EXT only had left+right,
# TOTAL only had top+bottom).
$bottomBorderRange = $ws.Range("A25:A26")
$bottomBorderRange.Borders.Item(7).LineStyle = 1
$bottomBorderRange.Borders.Item(10).LineStyle = 1
$bottomBorderRange.Borders.Item(9).LineStyle = 1
$ws.Range("A25").Borders.Item(8).LineStyle = -4142
$ws.Range("A26").Borders.Item(8).LineStyle = -4142

# Restore default top-left cell / selection (the saved view had scrolled
# to J16 with AB12 selected; the new save resets the scroll position and
# leaves Q17 selected).
[void]$ws.Range("A1").Select()
[void]$ws.Range("Q17").Select()
